$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6410857439041138
$ws.Range("B1").Value = 0.9297026991844177
$ws.Range("C1").Value = 1.098338007926941
$ws.Range("D1").Value = 3.91424822807312
$ws.Range("E1").Value = 2.327521085739136
